$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for this market/product subset.
# It belongs right after the existing row 9 (same cascading log order as
# the rest of the sheet), so insert a fresh row at position 10 — this
# pushes the former rows 10..101 down to 11..102, carrying their
# formatting (incl. the date-format style on column D) along with them.
$ws.Rows(10).Insert()

# Populate the newly inserted row 10 with the new record's data.
$ws.Cells.Item(10, 1).Value = 11
$ws.Cells.Item(10, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(10, 3).Value = "Bíobío"
$ws.Cells.Item(10, 4).Value = 44685
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100102
$ws.Cells.Item(10, 8).Value = "Cítricos"
$ws.Cells.Item(10, 9).Value = 100102004
$ws.Cells.Item(10, 10).Value = "Mandarina"
$ws.Cells.Item(10, 11).Value = "Murcott"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 100
$ws.Cells.Item(10, 14).Value = 10000
$ws.Cells.Item(10, 15).Value = 11000
$ws.Cells.Item(10, 16).Value = 10500
$ws.Cells.Item(10, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(10, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(10, 19).Value = 583
$ws.Cells.Item(10, 20).Value = 18
